# Weekly update: a new daily-price record was inserted for this market/product
# (row 30), which shifts every subsequent record down by one row
# (old row 30 -> new row 31, ..., old row 89 -> new row 90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30; this pushes rows 30..89 down to 31..90
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new record's data.
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44662
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 100112030
$ws.Range("G30").Value = "Poroto granado"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 460
$ws.Range("K30").Value = 24000
$ws.Range("L30").Value = 25000
$ws.Range("M30").Value = 24500
$ws.Range("N30").Value = "$/malla 25 kilos"
$ws.Range("O30").Value = "Provincia del Elquí"
$ws.Range("P30").Value = 980
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"
